# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (16-23) are reordered from newest-first
# (1807 .. 1712) to oldest-first (1712 .. 1807). Each period keeps its
# own "Valor Mora" (period 1807 = 9840, every other period = 29520),
# so after the reorder the 9840 value moves from row 16 down to row 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for rows 16..23
$periods = @("1712", "1801", "1802", "1803", "1804", "1805", "1806", "1807")

$startRow = 16
for ($idx = 0; $idx -lt $periods.Length; $idx++) {
    $row = $startRow + $idx
    $period = $periods[$idx]

    # Valor Mora follows the period itself: 1807 -> 9840, everything else -> 29520
    if ($period -eq "1807") {
        $valorMora = 9840
    } else {
        $valorMora = 29520
    }

    $ws.Cells.Item($row, 5).Value = $period
    $ws.Cells.Item($row, 6).Value = $valorMora
}
